$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 51 (volumeSize): SUNAT value "90Gi" -> "8Gi" ---
$ws.Range("C51").Value = "8Gi"

# --- Update existing row 52 (storageClassName): SUNAT value "storage-nfs" -> "nfs-storage" ---
$ws.Range("C52").Value = "nfs-storage"

# --- Append new variable comparison rows 53-77 (image registry/repo/version split into variables) ---
$ws.Cells.Item(53, 1).Value = "jenkinsImageRegistry"
$ws.Cells.Item(53, 2).Value = "docker.io"
$ws.Cells.Item(53, 3).Value = "vcf-np-w2-harbor-az1.sunat.peru"
$ws.Cells.Item(53, 2).Style = $ws.Cells.Item(50, 2).Style
$ws.Cells.Item(53, 3).Style = $ws.Cells.Item(50, 3).Style
$ws.Cells.Item(53, 4).Formula = '=$A$1&":"&A53&"|"&$B$1&":"&B53&"|"&$C$1&":"&C53'

$ws.Cells.Item(54, 1).Value = "jenkinsImageRepo"
$ws.Cells.Item(54, 2).Value = "bitnami/jenkins"
$ws.Cells.Item(54, 3).Value = "rayserve/jenkins"
$ws.Cells.Item(54, 2).Style = $ws.Cells.Item(50, 2).Style
$ws.Cells.Item(54, 3).Style = $ws.Cells.Item(50, 3).Style
$ws.Cells.Item(54, 4).Formula = '=$A$1&":"&A54&"|"&$B$1&":"&B54&"|"&$C$1&":"&C54'

$ws.Cells.Item(55, 1).Value = "jenkinsImageVersion"
$ws.Cells.Item(55, 2).Value = "2.401.1-debian-11-r0"
$ws.Cells.Item(55, 3).Value = "latest"
$ws.Cells.Item(55, 2).Style = $ws.Cells.Item(50, 2).Style
$ws.Cells.Item(55, 3).Style = $ws.Cells.Item(50, 3).Style
$ws.Cells.Item(55, 4).Formula = '=$A$1&":"&A55&"|"&$B$1&":"&B55&"|"&$C$1&":"&C55'

$ws.Cells.Item(56, 1).Value = "jenkinsAgentImageRegistry"
$ws.Cells.Item(56, 2).Value = "docker.io"
$ws.Cells.Item(56, 3).Value = "vcf-np-w2-harbor-az1.sunat.peru"
$ws.Cells.Item(56, 2).Style = $ws.Cells.Item(50, 2).Style
$ws.Cells.Item(56, 3).Style = $ws.Cells.Item(50, 3).Style
$ws.Cells.Item(56, 4).Formula = '=$A$1&":"&A56&"|"&$B$1&":"&B56&"|"&$C$1&":"&C56'

$ws.Cells.Item(57, 1).Value = "jenkinsAgentImageRepo"
$ws.Cells.Item(57, 2).Value = "bitnami/jenkins-agent"
$ws.Cells.Item(57, 3).Value = "rayserve/jenkins-agent"
$ws.Cells.Item(57, 2).Style = $ws.Cells.Item(50, 2).Style
$ws.Cells.Item(57, 3).Style = $ws.Cells.Item(50, 3).Style
$ws.Cells.Item(57, 4).Formula = '=$A$1&":"&A57&"|"&$B$1&":"&B57&"|"&$C$1&":"&C57'

$ws.Cells.Item(58, 1).Value = "jenkinsAgentImageVersion"
$ws.Cells.Item(58, 2).Value = "0.3107.0-debian-11-r36"
$ws.Cells.Item(58, 3).Value = "latest"
$ws.Cells.Item(58, 2).Style = $ws.Cells.Item(50, 2).Style
$ws.Cells.Item(58, 3).Style = $ws.Cells.Item(50, 3).Style
$ws.Cells.Item(58, 4).Formula = '=$A$1&":"&A58&"|"&$B$1&":"&B58&"|"&$C$1&":"&C58'

$ws.Cells.Item(59, 1).Value = "jenkinsInitContainerImage"
$ws.Cells.Item(59, 2).Value = "busybox"
$ws.Cells.Item(59, 3).Value = "vcf-np-w2-harbor-az1.sunat.peru/rayserve/busybox:latest"
$ws.Cells.Item(59, 2).Style = $ws.Cells.Item(50, 2).Style
$ws.Cells.Item(59, 3).Style = $ws.Cells.Item(50, 3).Style
$ws.Cells.Item(59, 4).Formula = '=$A$1&":"&A59&"|"&$B$1&":"&B59&"|"&$C$1&":"&C59'

$ws.Cells.Item(60, 1).Value = "jenkinsShellImageRegistry"
$ws.Cells.Item(60, 2).Value = "docker.io"
$ws.Cells.Item(60, 3).Value = "vcf-np-w2-harbor-az1.sunat.peru"
$ws.Cells.Item(60, 2).Style = $ws.Cells.Item(50, 2).Style
$ws.Cells.Item(60, 3).Style = $ws.Cells.Item(50, 3).Style
$ws.Cells.Item(60, 4).Formula = '=$A$1&":"&A60&"|"&$B$1&":"&B60&"|"&$C$1&":"&C60'

$ws.Cells.Item(61, 1).Value = "jenkinsShellImageRepo"
$ws.Cells.Item(61, 2).Value = "bitnami/bitnami-shell"
$ws.Cells.Item(61, 3).Value = "rayserve/bitnami-shell"
$ws.Cells.Item(61, 2).Style = $ws.Cells.Item(50, 2).Style
$ws.Cells.Item(61, 3).Style = $ws.Cells.Item(50, 3).Style
$ws.Cells.Item(61, 4).Formula = '=$A$1&":"&A61&"|"&$B$1&":"&B61&"|"&$C$1&":"&C61'

$ws.Cells.Item(62, 1).Value = "jenkinsShellImageVersion"
$ws.Cells.Item(62, 2).Value = "11-debian-11-r126"
$ws.Cells.Item(62, 3).Value = "11-r38"
$ws.Cells.Item(62, 2).Style = $ws.Cells.Item(50, 2).Style
$ws.Cells.Item(62, 3).Style = $ws.Cells.Item(50, 3).Style
$ws.Cells.Item(62, 4).Formula = '=$A$1&":"&A62&"|"&$B$1&":"&B62&"|"&$C$1&":"&C62'

$ws.Cells.Item(63, 1).Value = "minioImageRegistry"
$ws.Cells.Item(63, 2).Value = "docker.io"
$ws.Cells.Item(63, 3).Value = "vcf-np-w2-harbor-az1.sunat.peru"
$ws.Cells.Item(63, 2).Style = $ws.Cells.Item(50, 2).Style
$ws.Cells.Item(63, 3).Style = $ws.Cells.Item(50, 3).Style
$ws.Cells.Item(63, 4).Formula = '=$A$1&":"&A63&"|"&$B$1&":"&B63&"|"&$C$1&":"&C63'

$ws.Cells.Item(64, 1).Value = "minioImageRepo"
$ws.Cells.Item(64, 2).Value = "bitnami/minio"
$ws.Cells.Item(64, 3).Value = "rayserve/minio"
$ws.Cells.Item(64, 2).Style = $ws.Cells.Item(50, 2).Style
$ws.Cells.Item(64, 3).Style = $ws.Cells.Item(50, 3).Style
$ws.Cells.Item(64, 4).Formula = '=$A$1&":"&A64&"|"&$B$1&":"&B64&"|"&$C$1&":"&C64'

$ws.Cells.Item(65, 1).Value = "minioImageVersion"
$ws.Cells.Item(65, 2).Value = "2022.10.20-debian-11-r0"
$ws.Cells.Item(65, 3).Value = "latest"
$ws.Cells.Item(65, 2).Style = $ws.Cells.Item(50, 2).Style
$ws.Cells.Item(65, 3).Style = $ws.Cells.Item(50, 3).Style
$ws.Cells.Item(65, 4).Formula = '=$A$1&":"&A65&"|"&$B$1&":"&B65&"|"&$C$1&":"&C65'

$ws.Cells.Item(66, 1).Value = "minioClientImageRegistry"
$ws.Cells.Item(66, 2).Value = "docker.io"
$ws.Cells.Item(66, 3).Value = "vcf-np-w2-harbor-az1.sunat.peru"
$ws.Cells.Item(66, 2).Style = $ws.Cells.Item(50, 2).Style
$ws.Cells.Item(66, 3).Style = $ws.Cells.Item(50, 3).Style
$ws.Cells.Item(66, 4).Formula = '=$A$1&":"&A66&"|"&$B$1&":"&B66&"|"&$C$1&":"&C66'

$ws.Cells.Item(67, 1).Value = "minioClientImageRepo"
$ws.Cells.Item(67, 2).Value = "bitnami/minio-client"
$ws.Cells.Item(67, 3).Value = "rayserve/minio"
$ws.Cells.Item(67, 2).Style = $ws.Cells.Item(50, 2).Style
$ws.Cells.Item(67, 3).Style = $ws.Cells.Item(50, 3).Style
$ws.Cells.Item(67, 4).Formula = '=$A$1&":"&A67&"|"&$B$1&":"&B67&"|"&$C$1&":"&C67'

$ws.Cells.Item(68, 1).Value = "minioClientImageVersion"
$ws.Cells.Item(68, 2).Value = "2022.10.12-debian-11-r1"
$ws.Cells.Item(68, 3).Value = "2022.10.6-debian-11-r1"
$ws.Cells.Item(68, 2).Style = $ws.Cells.Item(50, 2).Style
$ws.Cells.Item(68, 3).Style = $ws.Cells.Item(50, 3).Style
$ws.Cells.Item(68, 4).Formula = '=$A$1&":"&A68&"|"&$B$1&":"&B68&"|"&$C$1&":"&C68'

$ws.Cells.Item(69, 1).Value = "minioShellImageRegistry"
$ws.Cells.Item(69, 2).Value = "docker.io"
$ws.Cells.Item(69, 3).Value = "vcf-np-w2-harbor-az1.sunat.peru"
$ws.Cells.Item(69, 2).Style = $ws.Cells.Item(50, 2).Style
$ws.Cells.Item(69, 3).Style = $ws.Cells.Item(50, 3).Style
$ws.Cells.Item(69, 4).Formula = '=$A$1&":"&A69&"|"&$B$1&":"&B69&"|"&$C$1&":"&C69'

$ws.Cells.Item(70, 1).Value = "minioShellImageRepo"
$ws.Cells.Item(70, 2).Value = "bitnami/bitnami-shell"
$ws.Cells.Item(70, 3).Value = "rayserve/minio"
$ws.Cells.Item(70, 2).Style = $ws.Cells.Item(50, 2).Style
$ws.Cells.Item(70, 3).Style = $ws.Cells.Item(50, 3).Style
$ws.Cells.Item(70, 4).Formula = '=$A$1&":"&A70&"|"&$B$1&":"&B70&"|"&$C$1&":"&C70'

$ws.Cells.Item(71, 1).Value = "minioShellImageVersion"
$ws.Cells.Item(71, 2).Value = "11-debian-11-r43"
$ws.Cells.Item(71, 3).Value = "latest"
$ws.Cells.Item(71, 2).Style = $ws.Cells.Item(50, 2).Style
$ws.Cells.Item(71, 3).Style = $ws.Cells.Item(50, 3).Style
$ws.Cells.Item(71, 4).Formula = '=$A$1&":"&A71&"|"&$B$1&":"&B71&"|"&$C$1&":"&C71'

$ws.Cells.Item(72, 1).Value = "postgresJupyterhubImageRegistry"
$ws.Cells.Item(72, 2).Value = "docker.io"
$ws.Cells.Item(72, 3).Value = "vcf-np-w2-harbor-az1.sunat.peru"
$ws.Cells.Item(72, 2).Style = $ws.Cells.Item(50, 2).Style
$ws.Cells.Item(72, 3).Style = $ws.Cells.Item(50, 3).Style
$ws.Cells.Item(72, 4).Formula = '=$A$1&":"&A72&"|"&$B$1&":"&B72&"|"&$C$1&":"&C72'

$ws.Cells.Item(73, 1).Value = "postgresJupyterhubImageRepo"
$ws.Cells.Item(73, 2).Value = "bitnami/postgresql"
$ws.Cells.Item(73, 3).Value = "rayserve/postgresql"
$ws.Cells.Item(73, 2).Style = $ws.Cells.Item(50, 2).Style
$ws.Cells.Item(73, 3).Style = $ws.Cells.Item(50, 3).Style
$ws.Cells.Item(73, 4).Formula = '=$A$1&":"&A73&"|"&$B$1&":"&B73&"|"&$C$1&":"&C73'

$ws.Cells.Item(74, 1).Value = "postgresJupyterhubImageVersion"
$ws.Cells.Item(74, 2).Value = "14.5.0-debian-11-r31"
$ws.Cells.Item(74, 3).Value = "latest"
$ws.Cells.Item(74, 2).Style = $ws.Cells.Item(50, 2).Style
$ws.Cells.Item(74, 3).Style = $ws.Cells.Item(50, 3).Style
$ws.Cells.Item(74, 4).Formula = '=$A$1&":"&A74&"|"&$B$1&":"&B74&"|"&$C$1&":"&C74'

$ws.Cells.Item(75, 1).Value = "postgresImageRegistry"
$ws.Cells.Item(75, 2).Value = "docker.io"
$ws.Cells.Item(75, 3).Value = "vcf-np-w2-harbor-az1.sunat.peru"
$ws.Cells.Item(75, 2).Style = $ws.Cells.Item(50, 2).Style
$ws.Cells.Item(75, 3).Style = $ws.Cells.Item(50, 3).Style
$ws.Cells.Item(75, 4).Formula = '=$A$1&":"&A75&"|"&$B$1&":"&B75&"|"&$C$1&":"&C75'

$ws.Cells.Item(76, 1).Value = "postgresImageRepo"
$ws.Cells.Item(76, 2).Value = "bitnami/postgresql"
$ws.Cells.Item(76, 3).Value = "rayserve/postgresql"
$ws.Cells.Item(76, 2).Style = $ws.Cells.Item(50, 2).Style
$ws.Cells.Item(76, 3).Style = $ws.Cells.Item(50, 3).Style
$ws.Cells.Item(76, 4).Formula = '=$A$1&":"&A76&"|"&$B$1&":"&B76&"|"&$C$1&":"&C76'

$ws.Cells.Item(77, 1).Value = "postgresImageVersion"
$ws.Cells.Item(77, 2).Value = "14.5.0-debian-11-r31"
$ws.Cells.Item(77, 3).Value = "latest"
$ws.Cells.Item(77, 2).Style = $ws.Cells.Item(50, 2).Style
$ws.Cells.Item(77, 3).Style = $ws.Cells.Item(50, 3).Style
$ws.Cells.Item(77, 4).Formula = '=$A$1&":"&A77&"|"&$B$1&":"&B77&"|"&$C$1&":"&C77'

# --- Update sheet view state ---
$ws.Application.ActiveWindow.ScrollRow = 46
$ws.Range("C52").Select()
